# Generate Report for handoff
#
# Refresh the "Latest Handoff Datetime" column (column D) on the per-locale
# report sheets. Rows whose last handoff datetime was still pointing at the
# previous handoff run ("...:08" for zh-cn / "...:26" for de-de, the time
# shared by every file that hasn't had a more recent, individual handoff)
# are bumped to reflect the new handoff timestamps recorded for this run.

$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$zhHandoffDatetime = "2016-02-16 10:49:08"
$deHandoffDatetime = "2016-02-16 10:49:24"

$rows = @(4, 6, 7, 8, 9, 10)

foreach ($r in $rows) {
    $wsZh.Range("D$r").Value = $zhHandoffDatetime
    $wsDe.Range("D$r").Value = $deHandoffDatetime
}

Write-Host "Updated Latest Handoff Datetime for rows: $rows"
